# Update the dSF column (column F) values for several rows to reflect
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    8  = 0
    12 = -2
    18 = 3
    20 = 4
    22 = 2
    24 = -1
    25 = -1
    30 = -3
    32 = 0
    34 = -2
    35 = -2
    37 = 9
    40 = 2
    47 = -2
    48 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
